$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038925902598975
$ws.Range("D2").Value = 1.045824259550994
$ws.Range("E2").Value = 1.037515753087569
$ws.Range("F2").Value = 1.053242771148646
$ws.Range("I2").Value = 1.035310397126852
$ws.Range("J2").Value = 1.044020694393578
$ws.Range("K2").Value = 1.048591416162014
$ws.Range("L2").Value = 1.040306415913151
$ws.Range("M2").Value = 1.05598927923025
$ws.Range("N2").Value = 1.018557384346191
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040146842485191
$ws.Range("D3").Value = 1.046953669195094
$ws.Range("E3").Value = 1.038559946249151
$ws.Range("F3").Value = 1.054536563024385
$ws.Range("I3").Value = 1.03554030401663
$ws.Range("J3").Value = 1.044885269418924
$ws.Range("K3").Value = 1.04953151742831
$ws.Range("L3").Value = 1.041159812397393
$ws.Range("M3").Value = 1.057094845379991
$ws.Range("N3").Value = 1.018846554919712
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040936008861296
$ws.Range("D4").Value = 1.047684007499962
$ws.Range("E4").Value = 1.039235228982136
$ws.Range("F4").Value = 1.055373541804755
$ws.Range("I4").Value = 1.035686628007292
$ws.Range("J4").Value = 1.045443347638329
$ws.Range("K4").Value = 1.050138779255029
$ws.Range("L4").Value = 1.041711039555703
$ws.Range("M4").Value = 1.057809480276827
$ws.Range("N4").Value = 1.019033130444385
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041267570210219
$ws.Range("D5").Value = 1.047990932057683
$ws.Range("E5").Value = 1.039519028216769
$ws.Range("F5").Value = 1.055725364159458
$ws.Range("I5").Value = 1.035747558436332
$ws.Range("J5").Value = 1.045677639645489
$ws.Range("K5").Value = 1.050393823122154
$ws.Range("L5").Value = 1.04194254273641
$ws.Range("M5").Value = 1.058109737438491
$ws.Range("N5").Value = 1.019111438615066
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041323228917765
$ws.Range("D6").Value = 1.048042459660218
$ws.Range("E6").Value = 1.039566674140602
$ws.Range("F6").Value = 1.055784434192803
$ws.Range("I6").Value = 1.035757754673691
$ws.Range("J6").Value = 1.045716959341459
$ws.Range("K6").Value = 1.050436631552449
$ws.Range("L6").Value = 1.041981399529479
$ws.Range("M6").Value = 1.058160141709118
$ws.Range("N6").Value = 1.019124579397076
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040940439999535
$ws.Range("D7").Value = 1.047688109069726
$ws.Range("E7").Value = 1.039239021471903
$ws.Range("F7").Value = 1.055378243040561
$ws.Range("I7").Value = 1.035687444457821
$ws.Range("J7").Value = 1.04544647953126
$ws.Range("K7").Value = 1.050142188140196
$ws.Range("L7").Value = 1.041714133825387
$ws.Range("M7").Value = 1.057813493014399
$ws.Range("N7").Value = 1.019034177305142
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03933870522569
$ws.Range("D8").Value = 1.046206046201736
$ws.Range("E8").Value = 1.037868722575609
$ws.Range("F8").Value = 1.053680055075482
$ws.Range("I8").Value = 1.03538860083847
$ws.Range("J8").Value = 1.044313163954667
$ws.Range("K8").Value = 1.048909345416478
$ws.Range("L8").Value = 1.04059502841125
$ws.Range("M8").Value = 1.056363064918164
$ws.Range("N8").Value = 1.018655222375349
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036509518583897
$ws.Range("D9").Value = 1.043590820687717
$ws.Range("E9").Value = 1.035451107840201
$ws.Range("F9").Value = 1.050686067013047
$ws.Range("I9").Value = 1.034843292285083
$ws.Range("J9").Value = 1.042305639025842
$ws.Range("K9").Value = 1.04672883345547
$ws.Range("L9").Value = 1.038615477214958
$ws.Range("M9").Value = 1.05380146962767
$ws.Range("N9").Value = 1.017983322042482
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034618711114529
$ws.Range("D10").Value = 1.041844775991766
$ws.Range("E10").Value = 1.033837275360043
$ws.Range("F10").Value = 1.048688879461536
$ws.Range("I10").Value = 1.034467156090151
$ws.Range("J10").Value = 1.040960150030895
$ws.Range("K10").Value = 1.045269621027673
$ws.Range("L10").Value = 1.037290618932799
$ws.Range("M10").Value = 1.05208974644373
$ws.Range("N10").Value = 1.017532581377436
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033798822914076
$ws.Range("D11").Value = 1.041088085738652
$ws.Range("E11").Value = 1.033137951316992
$ws.Range("F11").Value = 1.047823758517745
$ws.Range("I11").Value = 1.034301292032281
$ws.Range("J11").Value = 1.040375823158491
$ws.Range("K11").Value = 1.044636428837516
$ws.Range("L11").Value = 1.036715697492512
$ws.Range("M11").Value = 1.051347575851647
$ws.Range("N11").Value = 1.017336733261685
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03349410293628
$ws.Range("D12").Value = 1.040806918550213
$ws.Range("E12").Value = 1.032878110571421
$ws.Range("F12").Value = 1.047502362195293
$ws.Range("I12").Value = 1.034239232290016
$ws.Range("J12").Value = 1.040158517523385
$ws.Range("K12").Value = 1.044401029054725
$ws.Range("L12").Value = 1.036501956587866
$ws.Range("M12").Value = 1.05107175036612
$ws.Range("N12").Value = 1.017263884610344
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033559474468634
$ws.Range("D13").Value = 1.040867234373036
$ws.Range("E13").Value = 1.032933850974807
$ws.Range("F13").Value = 1.047571305193344
$ws.Range("I13").Value = 1.034252564707878
$ws.Range("J13").Value = 1.040205142145902
$ws.Range("K13").Value = 1.044451532373915
$ws.Range("L13").Value = 1.036547813341522
$ws.Range("M13").Value = 1.05113092271125
$ws.Range("N13").Value = 1.017279515516343
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.0337736382955
$ws.Range("D14").Value = 1.041064846399177
$ws.Range("E14").Value = 1.033116474456637
$ws.Range("F14").Value = 1.047797192873183
$ws.Range("I14").Value = 1.034296171346773
$ws.Range("J14").Value = 1.040357865935174
$ws.Range("K14").Value = 1.044616974792319
$ws.Range("L14").Value = 1.036698033478693
$ws.Range("M14").Value = 1.051324779111658
$ws.Range("N14").Value = 1.017330713652589
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033905568144719
$ws.Range("D15").Value = 1.041186588535425
$ws.Range("E15").Value = 1.033228984064568
$ws.Range("F15").Value = 1.047936362763422
$ws.Range("I15").Value = 1.034322979131279
$ws.Range("J15").Value = 1.040451929507175
$ws.Range("K15").Value = 1.044718882228452
$ws.Range("L15").Value = 1.036790563899753
$ws.Range("M15").Value = 1.05144420044674
$ws.Range("N15").Value = 1.01736224498739
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034673099755137
$ws.Range("D16").Value = 1.041894981316299
$ws.Range("E16").Value = 1.033883675989561
$ws.Range("F16").Value = 1.048746287577044
$ws.Range("I16").Value = 1.034478100802755
$ws.Range("J16").Value = 1.040998893444302
$ws.Range("K16").Value = 1.045311615417019
$ws.Range("L16").Value = 1.037328748103636
$ws.Range("M16").Value = 1.052138980928444
$ws.Range("N16").Value = 1.017545564909996
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03515424021243
$ws.Range("D17").Value = 1.042339164005641
$ws.Range("E17").Value = 1.034294205300164
$ws.Range("F17").Value = 1.04925424276211
$ws.Range("I17").Value = 1.034574602411513
$ws.Range("J17").Value = 1.041341527223458
$ws.Range("K17").Value = 1.045683060058359
$ws.Range("L17").Value = 1.037666001143507
$ws.Range("M17").Value = 1.052574533338801
$ws.Range("N17").Value = 1.017660375657077
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035434769980931
$ws.Range("D18").Value = 1.042598186532733
$ws.Range("E18").Value = 1.034533609683487
$ws.Range("F18").Value = 1.049550493219753
$ws.Range("I18").Value = 1.034630601269468
$ws.Range("J18").Value = 1.041541213771993
$ws.Range("K18").Value = 1.045899587945031
$ws.Range("L18").Value = 1.037862594744925
$ws.Range("M18").Value = 1.052828489006619
$ws.Range("N18").Value = 1.017727277792674
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035530404599777
$ws.Range("D19").Value = 1.042686496122921
$ws.Range("E19").Value = 1.034615231897701
$ws.Range("F19").Value = 1.049651501765695
$ws.Range("I19").Value = 1.034649646446689
$ws.Range("J19").Value = 1.041609273639013
$ws.Range("K19").Value = 1.045973396459714
$ws.Range("L19").Value = 1.037929607716727
$ws.Range("M19").Value = 1.052915065343516
$ws.Range("N19").Value = 1.017750078680014
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035102629934382
$ws.Range("D20").Value = 1.042291513809143
$ws.Range("E20").Value = 1.034250164620383
$ws.Range("F20").Value = 1.04919974722875
$ws.Range("I20").Value = 1.034564278588562
$ws.Range("J20").Value = 1.041304783040903
$ws.Range("K20").Value = 1.045643220973814
$ws.Range("L20").Value = 1.037629829550763
$ws.Range("M20").Value = 1.052527812500586
$ws.Range("N20").Value = 1.017648064281678
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033710577298969
$ws.Range("D21").Value = 1.041006657327651
$ws.Range("E21").Value = 1.033062698635298
$ws.Range("F21").Value = 1.047730676037836
$ws.Range("I21").Value = 1.034283342717421
$ws.Range("J21").Value = 1.040312899796408
$ws.Range("K21").Value = 1.044568261772727
$ws.Range("L21").Value = 1.036653802637865
$ws.Range("M21").Value = 1.0512676974028
$ws.Range("N21").Value = 1.01731563989659
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032834312865644
$ws.Range("D22").Value = 1.040198244607969
$ws.Range("E22").Value = 1.032315623996666
$ws.Range("F22").Value = 1.046806710996121
$ws.Range("I22").Value = 1.034104100764015
$ws.Range("J22").Value = 1.039687754565267
$ws.Range("K22").Value = 1.043891210834904
$ws.Range("L22").Value = 1.036039038568833
$ws.Range("M22").Value = 1.050474542247968
$ws.Range("N22").Value = 1.017106041382966
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033298935555574
$ws.Range("D23").Value = 1.040626854548351
$ws.Range("E23").Value = 1.032711707379861
$ws.Range("F23").Value = 1.047296551941715
$ws.Range("I23").Value = 1.03419936755819
$ws.Range("J23").Value = 1.040019299583457
$ws.Range("K23").Value = 1.044250241057285
$ws.Range("L23").Value = 1.03636504125156
$ws.Range("M23").Value = 1.050895092208162
$ws.Range("N23").Value = 1.017217209663561
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035125950741368
$ws.Range("D24").Value = 1.042313045072566
$ws.Range("E24").Value = 1.034270064863678
$ws.Range("F24").Value = 1.049224371510931
$ws.Range("I24").Value = 1.034568944372448
$ws.Range("J24").Value = 1.041321386668073
$ws.Range("K24").Value = 1.04566122294135
$ws.Range("L24").Value = 1.037646174308527
$ws.Range("M24").Value = 1.052548923928829
$ws.Range("N24").Value = 1.017653627463224
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037241744654515
$ws.Range("D25").Value = 1.044267362222717
$ws.Range("E25").Value = 1.036076480752786
$ws.Range("F25").Value = 1.051460284771143
$ws.Range("I25").Value = 1.034986486498896
$ws.Range("J25").Value = 1.04282588366386
$ws.Range("K25").Value = 1.047293516620539
$ws.Range("L25").Value = 1.039128141389411
$ws.Range("M25").Value = 1.054464398093522
$ws.Range("N25").Value = 1.018157517071797
